$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2023-12-11 00:18:56"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0.001
$ws.Range("I7").Value = 0.01
$ws.Range("J7").Value = 0.003
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 512
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 0.9230769230769231

$ws.Range("A8").Value = "2023-12-12 19:20:40"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 13
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.001
$ws.Range("I8").Value = 0.01
$ws.Range("J8").Value = 0.003
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 512
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 7
$ws.Range("O8").Value = 1.625
